$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Apoe","Lrp8","ECs",3,1,67.77251700000001,203.317551,0.4079637943863715,0.4079637943863715,3,1,3.211751333333333,9.635254,0.9274105493513785,0.9274105493513782,217.668471838106,1959.016246542954,0.3783499266673376,0.3783499266673375),
    @("ECs","Apoe","Lrp8","FAPs",3,1,67.77251700000001,203.317551,0.4079637943863715,0.4079637943863715,2,0.6666666666666666,0.2513873333333334,0.754162,0.07258945064862164,0.07258945064862163,17.037152321918,153.334370897262,0.02961386771903394,0.02961386771903394),
    @("FAPs","Apoe","Lrp8","ECs",3,1,60.97760633333333,182.932819,0.3670611149405164,0.3670611149405164,3,1,3.211751333333333,9.635254,0.9274105493513785,0.9274105493513782,195.8449084445585,1762.604176001026,0.3404163502525138,0.3404163502525137),
    @("FAPs","Apoe","Lrp8","FAPs",3,1,60.97760633333333,182.932819,0.3670611149405164,0.3670611149405164,2,0.6666666666666666,0.2513873333333334,0.754162,0.07258945064862164,0.07258945064862163,15.32899784918644,137.960980642678,0.02664476468800265,0.02664476468800264),
    @("sCs","Apoe","Lrp8","ECs",3,1,37.37372866666667,112.121186,0.2249750906731122,0.2249750906731122,3,1,3.211751333333333,9.635254,0.9274105493513785,0.9274105493513782,120.0351228768049,1080.316105891244,0.2086442724315272,0.2086442724315271),
    @("sCs","Apoe","Lrp8","FAPs",3,1,37.37372866666667,112.121186,0.2249750906731122,0.2249750906731122,2,0.6666666666666666,0.2513873333333334,0.754162,0.07258945064862164,0.07258945064862163,9.39528198623689,84.557537876132,0.01633081824158506,0.01633081824158505)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
